# Add new GPU model rows (RX 7814 XT .. RX 7823 XT) and extra FPS columns to the CPU sheet,
# mirroring the existing block pattern (B:D, E:G, H:J, K:M, N, O:Q, R:T) through row 16.
$wb = $excel.ActiveWorkbook
$wsCpu = $wb.Worksheets.Item("CPU")

$wsCpu.Cells.Item(1, 1).Value = 'GPU 名稱'
$wsCpu.Cells.Item(1, 2).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 3).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 4).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 5).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 6).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 7).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 8).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 9).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 10).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 11).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 12).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 13).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 14).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 15).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 16).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 17).Value = 'FPS (4K)'
$wsCpu.Cells.Item(1, 18).Value = 'FPS (1080p)'
$wsCpu.Cells.Item(1, 19).Value = 'FPS (1440p)'
$wsCpu.Cells.Item(1, 20).Value = 'FPS (4K)'
$wsCpu.Cells.Item(2, 1).Value = 'RTX 4070'
$wsCpu.Cells.Item(2, 2).Value = 120
$wsCpu.Cells.Item(2, 3).Value = 95
$wsCpu.Cells.Item(2, 4).Value = 60
$wsCpu.Cells.Item(2, 5).Value = 120
$wsCpu.Cells.Item(2, 6).Value = 95
$wsCpu.Cells.Item(2, 7).Value = 60
$wsCpu.Cells.Item(2, 8).Value = 120
$wsCpu.Cells.Item(2, 9).Value = 95
$wsCpu.Cells.Item(2, 10).Value = 60
$wsCpu.Cells.Item(2, 11).Value = 120
$wsCpu.Cells.Item(2, 12).Value = 95
$wsCpu.Cells.Item(2, 13).Value = 60
$wsCpu.Cells.Item(2, 14).Value = 60
$wsCpu.Cells.Item(2, 15).Value = 120
$wsCpu.Cells.Item(2, 16).Value = 95
$wsCpu.Cells.Item(2, 17).Value = 60
$wsCpu.Cells.Item(2, 18).Value = 120
$wsCpu.Cells.Item(2, 19).Value = 95
$wsCpu.Cells.Item(2, 20).Value = 60
$wsCpu.Cells.Item(3, 1).Value = 'RX 7800 XT'
$wsCpu.Cells.Item(3, 2).Value = 110
$wsCpu.Cells.Item(3, 3).Value = 90
$wsCpu.Cells.Item(3, 4).Value = 55
$wsCpu.Cells.Item(3, 5).Value = 110
$wsCpu.Cells.Item(3, 6).Value = 90
$wsCpu.Cells.Item(3, 7).Value = 55
$wsCpu.Cells.Item(3, 8).Value = 110
$wsCpu.Cells.Item(3, 9).Value = 90
$wsCpu.Cells.Item(3, 10).Value = 55
$wsCpu.Cells.Item(3, 11).Value = 110
$wsCpu.Cells.Item(3, 12).Value = 90
$wsCpu.Cells.Item(3, 13).Value = 55
$wsCpu.Cells.Item(3, 14).Value = 55
$wsCpu.Cells.Item(3, 15).Value = 110
$wsCpu.Cells.Item(3, 16).Value = 90
$wsCpu.Cells.Item(3, 17).Value = 55
$wsCpu.Cells.Item(3, 18).Value = 110
$wsCpu.Cells.Item(3, 19).Value = 90
$wsCpu.Cells.Item(3, 20).Value = 55
$wsCpu.Cells.Item(4, 1).Value = 'RX 7801 XT'
$wsCpu.Cells.Item(4, 2).Value = 111
$wsCpu.Cells.Item(4, 3).Value = 91
$wsCpu.Cells.Item(4, 4).Value = 56
$wsCpu.Cells.Item(4, 5).Value = 111
$wsCpu.Cells.Item(4, 6).Value = 91
$wsCpu.Cells.Item(4, 7).Value = 56
$wsCpu.Cells.Item(4, 8).Value = 111
$wsCpu.Cells.Item(4, 9).Value = 91
$wsCpu.Cells.Item(4, 10).Value = 56
$wsCpu.Cells.Item(4, 11).Value = 111
$wsCpu.Cells.Item(4, 12).Value = 91
$wsCpu.Cells.Item(4, 13).Value = 56
$wsCpu.Cells.Item(4, 14).Value = 56
$wsCpu.Cells.Item(4, 15).Value = 111
$wsCpu.Cells.Item(4, 16).Value = 91
$wsCpu.Cells.Item(4, 17).Value = 56
$wsCpu.Cells.Item(4, 18).Value = 111
$wsCpu.Cells.Item(4, 19).Value = 91
$wsCpu.Cells.Item(4, 20).Value = 56
$wsCpu.Cells.Item(5, 1).Value = 'RX 7802 XT'
$wsCpu.Cells.Item(5, 2).Value = 112
$wsCpu.Cells.Item(5, 3).Value = 92
$wsCpu.Cells.Item(5, 4).Value = 57
$wsCpu.Cells.Item(5, 5).Value = 112
$wsCpu.Cells.Item(5, 6).Value = 92
$wsCpu.Cells.Item(5, 7).Value = 57
$wsCpu.Cells.Item(5, 8).Value = 112
$wsCpu.Cells.Item(5, 9).Value = 92
$wsCpu.Cells.Item(5, 10).Value = 57
$wsCpu.Cells.Item(5, 11).Value = 112
$wsCpu.Cells.Item(5, 12).Value = 92
$wsCpu.Cells.Item(5, 13).Value = 57
$wsCpu.Cells.Item(5, 14).Value = 57
$wsCpu.Cells.Item(5, 15).Value = 112
$wsCpu.Cells.Item(5, 16).Value = 92
$wsCpu.Cells.Item(5, 17).Value = 57
$wsCpu.Cells.Item(5, 18).Value = 112
$wsCpu.Cells.Item(5, 19).Value = 92
$wsCpu.Cells.Item(5, 20).Value = 57
$wsCpu.Cells.Item(6, 1).Value = 'RX 7813 XT'
$wsCpu.Cells.Item(6, 2).Value = 123
$wsCpu.Cells.Item(6, 3).Value = 103
$wsCpu.Cells.Item(6, 4).Value = 68
$wsCpu.Cells.Item(6, 5).Value = 113
$wsCpu.Cells.Item(6, 6).Value = 93
$wsCpu.Cells.Item(6, 7).Value = 58
$wsCpu.Cells.Item(6, 8).Value = 113
$wsCpu.Cells.Item(6, 9).Value = 93
$wsCpu.Cells.Item(6, 10).Value = 58
$wsCpu.Cells.Item(6, 11).Value = 113
$wsCpu.Cells.Item(6, 12).Value = 93
$wsCpu.Cells.Item(6, 13).Value = 58
$wsCpu.Cells.Item(6, 14).Value = 58
$wsCpu.Cells.Item(6, 15).Value = 113
$wsCpu.Cells.Item(6, 16).Value = 93
$wsCpu.Cells.Item(6, 17).Value = 58
$wsCpu.Cells.Item(6, 18).Value = 113
$wsCpu.Cells.Item(6, 19).Value = 93
$wsCpu.Cells.Item(6, 20).Value = 58
$wsCpu.Cells.Item(7, 1).Value = 'RX 7814 XT'
$wsCpu.Cells.Item(7, 2).Value = 124
$wsCpu.Cells.Item(7, 3).Value = 103
$wsCpu.Cells.Item(7, 4).Value = 68
$wsCpu.Cells.Item(7, 5).Value = 114
$wsCpu.Cells.Item(7, 6).Value = 94
$wsCpu.Cells.Item(7, 7).Value = 59
$wsCpu.Cells.Item(7, 8).Value = 114
$wsCpu.Cells.Item(7, 9).Value = 94
$wsCpu.Cells.Item(7, 10).Value = 59
$wsCpu.Cells.Item(7, 11).Value = 114
$wsCpu.Cells.Item(7, 12).Value = 94
$wsCpu.Cells.Item(7, 13).Value = 59
$wsCpu.Cells.Item(7, 14).Value = 59
$wsCpu.Cells.Item(7, 15).Value = 114
$wsCpu.Cells.Item(7, 16).Value = 94
$wsCpu.Cells.Item(7, 17).Value = 59
$wsCpu.Cells.Item(7, 18).Value = 114
$wsCpu.Cells.Item(7, 19).Value = 94
$wsCpu.Cells.Item(7, 20).Value = 59
$wsCpu.Cells.Item(8, 1).Value = 'RX 7815 XT'
$wsCpu.Cells.Item(8, 2).Value = 125
$wsCpu.Cells.Item(8, 3).Value = 103
$wsCpu.Cells.Item(8, 4).Value = 68
$wsCpu.Cells.Item(8, 5).Value = 115
$wsCpu.Cells.Item(8, 6).Value = 95
$wsCpu.Cells.Item(8, 7).Value = 60
$wsCpu.Cells.Item(8, 8).Value = 115
$wsCpu.Cells.Item(8, 9).Value = 95
$wsCpu.Cells.Item(8, 10).Value = 60
$wsCpu.Cells.Item(8, 11).Value = 115
$wsCpu.Cells.Item(8, 12).Value = 95
$wsCpu.Cells.Item(8, 13).Value = 60
$wsCpu.Cells.Item(8, 14).Value = 60
$wsCpu.Cells.Item(8, 15).Value = 115
$wsCpu.Cells.Item(8, 16).Value = 95
$wsCpu.Cells.Item(8, 17).Value = 60
$wsCpu.Cells.Item(8, 18).Value = 115
$wsCpu.Cells.Item(8, 19).Value = 95
$wsCpu.Cells.Item(8, 20).Value = 60
$wsCpu.Cells.Item(9, 1).Value = 'RX 7816 XT'
$wsCpu.Cells.Item(9, 2).Value = 126
$wsCpu.Cells.Item(9, 3).Value = 103
$wsCpu.Cells.Item(9, 4).Value = 68
$wsCpu.Cells.Item(9, 5).Value = 116
$wsCpu.Cells.Item(9, 6).Value = 96
$wsCpu.Cells.Item(9, 7).Value = 61
$wsCpu.Cells.Item(9, 8).Value = 116
$wsCpu.Cells.Item(9, 9).Value = 96
$wsCpu.Cells.Item(9, 10).Value = 61
$wsCpu.Cells.Item(9, 11).Value = 116
$wsCpu.Cells.Item(9, 12).Value = 96
$wsCpu.Cells.Item(9, 13).Value = 61
$wsCpu.Cells.Item(9, 14).Value = 61
$wsCpu.Cells.Item(9, 15).Value = 116
$wsCpu.Cells.Item(9, 16).Value = 96
$wsCpu.Cells.Item(9, 17).Value = 61
$wsCpu.Cells.Item(9, 18).Value = 116
$wsCpu.Cells.Item(9, 19).Value = 96
$wsCpu.Cells.Item(9, 20).Value = 61
$wsCpu.Cells.Item(10, 1).Value = 'RX 7817 XT'
$wsCpu.Cells.Item(10, 2).Value = 127
$wsCpu.Cells.Item(10, 3).Value = 103
$wsCpu.Cells.Item(10, 4).Value = 68
$wsCpu.Cells.Item(10, 5).Value = 117
$wsCpu.Cells.Item(10, 6).Value = 97
$wsCpu.Cells.Item(10, 7).Value = 62
$wsCpu.Cells.Item(10, 8).Value = 117
$wsCpu.Cells.Item(10, 9).Value = 97
$wsCpu.Cells.Item(10, 10).Value = 62
$wsCpu.Cells.Item(10, 11).Value = 117
$wsCpu.Cells.Item(10, 12).Value = 97
$wsCpu.Cells.Item(10, 13).Value = 62
$wsCpu.Cells.Item(10, 14).Value = 62
$wsCpu.Cells.Item(10, 15).Value = 117
$wsCpu.Cells.Item(10, 16).Value = 97
$wsCpu.Cells.Item(10, 17).Value = 62
$wsCpu.Cells.Item(10, 18).Value = 117
$wsCpu.Cells.Item(10, 19).Value = 97
$wsCpu.Cells.Item(10, 20).Value = 62
$wsCpu.Cells.Item(11, 1).Value = 'RX 7818 XT'
$wsCpu.Cells.Item(11, 2).Value = 128
$wsCpu.Cells.Item(11, 3).Value = 103
$wsCpu.Cells.Item(11, 4).Value = 68
$wsCpu.Cells.Item(11, 5).Value = 118
$wsCpu.Cells.Item(11, 6).Value = 98
$wsCpu.Cells.Item(11, 7).Value = 63
$wsCpu.Cells.Item(11, 8).Value = 118
$wsCpu.Cells.Item(11, 9).Value = 98
$wsCpu.Cells.Item(11, 10).Value = 63
$wsCpu.Cells.Item(11, 11).Value = 118
$wsCpu.Cells.Item(11, 12).Value = 98
$wsCpu.Cells.Item(11, 13).Value = 63
$wsCpu.Cells.Item(11, 14).Value = 63
$wsCpu.Cells.Item(11, 15).Value = 118
$wsCpu.Cells.Item(11, 16).Value = 98
$wsCpu.Cells.Item(11, 17).Value = 63
$wsCpu.Cells.Item(11, 18).Value = 118
$wsCpu.Cells.Item(11, 19).Value = 98
$wsCpu.Cells.Item(11, 20).Value = 63
$wsCpu.Cells.Item(12, 1).Value = 'RX 7819 XT'
$wsCpu.Cells.Item(12, 2).Value = 129
$wsCpu.Cells.Item(12, 3).Value = 103
$wsCpu.Cells.Item(12, 4).Value = 68
$wsCpu.Cells.Item(12, 5).Value = 119
$wsCpu.Cells.Item(12, 6).Value = 99
$wsCpu.Cells.Item(12, 7).Value = 64
$wsCpu.Cells.Item(12, 8).Value = 119
$wsCpu.Cells.Item(12, 9).Value = 99
$wsCpu.Cells.Item(12, 10).Value = 64
$wsCpu.Cells.Item(12, 11).Value = 119
$wsCpu.Cells.Item(12, 12).Value = 99
$wsCpu.Cells.Item(12, 13).Value = 64
$wsCpu.Cells.Item(12, 14).Value = 64
$wsCpu.Cells.Item(12, 15).Value = 119
$wsCpu.Cells.Item(12, 16).Value = 99
$wsCpu.Cells.Item(12, 17).Value = 64
$wsCpu.Cells.Item(12, 18).Value = 119
$wsCpu.Cells.Item(12, 19).Value = 99
$wsCpu.Cells.Item(12, 20).Value = 64
$wsCpu.Cells.Item(13, 1).Value = 'RX 7820 XT'
$wsCpu.Cells.Item(13, 2).Value = 130
$wsCpu.Cells.Item(13, 3).Value = 103
$wsCpu.Cells.Item(13, 4).Value = 68
$wsCpu.Cells.Item(13, 5).Value = 120
$wsCpu.Cells.Item(13, 6).Value = 100
$wsCpu.Cells.Item(13, 7).Value = 65
$wsCpu.Cells.Item(13, 8).Value = 120
$wsCpu.Cells.Item(13, 9).Value = 100
$wsCpu.Cells.Item(13, 10).Value = 65
$wsCpu.Cells.Item(13, 11).Value = 120
$wsCpu.Cells.Item(13, 12).Value = 100
$wsCpu.Cells.Item(13, 13).Value = 65
$wsCpu.Cells.Item(13, 14).Value = 65
$wsCpu.Cells.Item(13, 15).Value = 120
$wsCpu.Cells.Item(13, 16).Value = 100
$wsCpu.Cells.Item(13, 17).Value = 65
$wsCpu.Cells.Item(13, 18).Value = 120
$wsCpu.Cells.Item(13, 19).Value = 100
$wsCpu.Cells.Item(13, 20).Value = 65
$wsCpu.Cells.Item(14, 1).Value = 'RX 7821 XT'
$wsCpu.Cells.Item(14, 2).Value = 131
$wsCpu.Cells.Item(14, 3).Value = 103
$wsCpu.Cells.Item(14, 4).Value = 68
$wsCpu.Cells.Item(14, 5).Value = 121
$wsCpu.Cells.Item(14, 6).Value = 101
$wsCpu.Cells.Item(14, 7).Value = 66
$wsCpu.Cells.Item(14, 8).Value = 121
$wsCpu.Cells.Item(14, 9).Value = 101
$wsCpu.Cells.Item(14, 10).Value = 66
$wsCpu.Cells.Item(14, 11).Value = 121
$wsCpu.Cells.Item(14, 12).Value = 101
$wsCpu.Cells.Item(14, 13).Value = 66
$wsCpu.Cells.Item(14, 14).Value = 66
$wsCpu.Cells.Item(14, 15).Value = 121
$wsCpu.Cells.Item(14, 16).Value = 101
$wsCpu.Cells.Item(14, 17).Value = 66
$wsCpu.Cells.Item(14, 18).Value = 121
$wsCpu.Cells.Item(14, 19).Value = 101
$wsCpu.Cells.Item(14, 20).Value = 66
$wsCpu.Cells.Item(15, 1).Value = 'RX 7822 XT'
$wsCpu.Cells.Item(15, 2).Value = 132
$wsCpu.Cells.Item(15, 3).Value = 103
$wsCpu.Cells.Item(15, 4).Value = 68
$wsCpu.Cells.Item(15, 5).Value = 122
$wsCpu.Cells.Item(15, 6).Value = 102
$wsCpu.Cells.Item(15, 7).Value = 67
$wsCpu.Cells.Item(15, 8).Value = 122
$wsCpu.Cells.Item(15, 9).Value = 102
$wsCpu.Cells.Item(15, 10).Value = 67
$wsCpu.Cells.Item(15, 11).Value = 122
$wsCpu.Cells.Item(15, 12).Value = 102
$wsCpu.Cells.Item(15, 13).Value = 67
$wsCpu.Cells.Item(15, 14).Value = 67
$wsCpu.Cells.Item(15, 15).Value = 122
$wsCpu.Cells.Item(15, 16).Value = 102
$wsCpu.Cells.Item(15, 17).Value = 67
$wsCpu.Cells.Item(15, 18).Value = 122
$wsCpu.Cells.Item(15, 19).Value = 102
$wsCpu.Cells.Item(15, 20).Value = 67
$wsCpu.Cells.Item(16, 1).Value = 'RX 7823 XT'
$wsCpu.Cells.Item(16, 2).Value = 133
$wsCpu.Cells.Item(16, 3).Value = 103
$wsCpu.Cells.Item(16, 4).Value = 68
$wsCpu.Cells.Item(16, 5).Value = 123
$wsCpu.Cells.Item(16, 6).Value = 103
$wsCpu.Cells.Item(16, 7).Value = 68
$wsCpu.Cells.Item(16, 8).Value = 123
$wsCpu.Cells.Item(16, 9).Value = 103
$wsCpu.Cells.Item(16, 10).Value = 68
$wsCpu.Cells.Item(16, 11).Value = 123
$wsCpu.Cells.Item(16, 12).Value = 103
$wsCpu.Cells.Item(16, 13).Value = 68
$wsCpu.Cells.Item(16, 14).Value = 68
$wsCpu.Cells.Item(16, 15).Value = 123
$wsCpu.Cells.Item(16, 16).Value = 103
$wsCpu.Cells.Item(16, 17).Value = 68
$wsCpu.Cells.Item(16, 18).Value = 123
$wsCpu.Cells.Item(16, 19).Value = 103
$wsCpu.Cells.Item(16, 20).Value = 68

# Restore view selections recorded in the saved workbook.
$wsSpec = $wb.Worksheets.Item("SPEC")
$null = $wsSpec.Range("B1:D16").Select()

$null = $wsCpu.Range("I14").Select()
